# DevSecOpsPipeline slide update:
#  - add security-tool call-out labels (new paragraphs) under several
#    pipeline-stage boxes
#  - nudge/resize a handful of boxes + connectors to make room for the
#    new second line of text
#
# NOTE on precision: PowerPoint's Shape.Left/Top/Width/Height are expressed
# in points (1 pt = 12700 EMU) and are stored internally as single-precision
# floats, so an EMU value round-tripped through `Width = emu/12700.0` can
# land 1 EMU away from the intended integer. The literal point constants
# below were solved so that, after the engine's float32 quantization, they
# reproduce the exact target EMU extents/offsets from the source diff.

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1. Connector (id 60) fed from "Infrastructure Code Vulnerability..."
#    box up to "Container Vulnerability Scanning" box shrinks in width
#    now that the box it anchors to changed: ext cx 834988 -> 653999
# ---------------------------------------------------------------------
$sh60 = Get-ShapeById $s.Shapes 60
$sh60.Width = 51.49598693847656   # -> 653999 EMU (was 834988 EMU)

# ---------------------------------------------------------------------
# 2. "Static code analysis" box (id 11) gains a "[SonarQube]" line
# ---------------------------------------------------------------------
$sh11 = Get-ShapeById $s.Shapes 11
$sh11.TextFrame.TextRange.InsertAfter("`r[SonarQube]") | Out-Null

# ---------------------------------------------------------------------
# 3. "Library Vulnerability Scanning" box (id 14) gains a
#    "[Jfrog XRay]" line (built run-by-run to match source run breaks)
# ---------------------------------------------------------------------
$sh14 = Get-ShapeById $s.Shapes 14
$tr14 = $sh14.TextFrame.TextRange
$tr14.InsertAfter("`r[J") | Out-Null
$tr14.InsertAfter("f") | Out-Null
$tr14.InsertAfter("rog ") | Out-Null
$tr14.InsertAfter("XRay") | Out-Null
$tr14.InsertAfter("]") | Out-Null

# ---------------------------------------------------------------------
# 4. "Container Vulnerability Scanning" box (id 21) gains a
#    "[Clair by CoreOs]" line
# ---------------------------------------------------------------------
$sh21 = Get-ShapeById $s.Shapes 21
$tr21 = $sh21.TextFrame.TextRange
$tr21.InsertAfter("`r[Clair by ") | Out-Null
$tr21.InsertAfter("CoreOs") | Out-Null
$tr21.InsertAfter("]") | Out-Null

# ---------------------------------------------------------------------
# 5. "Infrastructure Code Vulnerability and Compiance Scanning" box
#    (id 26) widens and shifts left to fit the longer text, and gains
#    a "[Terrascan]" line
# ---------------------------------------------------------------------
$sh26 = Get-ShapeById $s.Shapes 26
$sh26.Left = 177.98220825195312    # -> 2260374 EMU (was 2441363 EMU)
$sh26.Width = 127.29772186279297   # -> 1616681 EMU (was 1308000 EMU)
$tr26 = $sh26.TextFrame.TextRange
$tr26.InsertAfter("`r[") | Out-Null
$tr26.InsertAfter("Terrascan") | Out-Null
$tr26.InsertAfter("]") | Out-Null

# ---------------------------------------------------------------------
# 6. Connector (id 44) feeding into box id 26 is repositioned/resized
#    to keep connecting correctly after box 26 moved
# ---------------------------------------------------------------------
$sh44 = Get-ShapeById $s.Shapes 44
$sh44.Left = 137.17276000976562    # -> 1742094 EMU (was 1755418 EMU)
$sh44.Top = 195.28985595703125     # -> 2480181 EMU (was 2466857 EMU)
$sh44.Height = 166.6409454345703   # -> 2116340 EMU (was 2142988 EMU)

# ---------------------------------------------------------------------
# 7. Connector (id 47) feeding into box id 26 from the right is
#    repositioned/resized likewise
# ---------------------------------------------------------------------
$sh47 = Get-ShapeById $s.Shapes 47
$sh47.Left = 305.2799377441406     # -> 3877055 EMU (was 3749363 EMU)
$sh47.Width = 110.25260162353516   # -> 1400208 EMU (was 1527900 EMU)

# ---------------------------------------------------------------------
# 8. "Compliance Scanning" box (id 50) gains a "[Amazon Inspector]" line
# ---------------------------------------------------------------------
$sh50 = Get-ShapeById $s.Shapes 50
$sh50.TextFrame.TextRange.InsertAfter("`r[Amazon Inspector]") | Out-Null
